$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update header row (row 1): shift values in C1, D1, E1
$ws.Range("C1").Value = "prediction"
$ws.Range("D1").Value = "rejection-f"
$ws.Range("E1").Value = "max"

# Update data rows 2-7: column C becomes the species string,
# column D stays the species string, column E becomes 0.5
for ($r = 2; $r -le 7; $r++) {
    $ws.Cells.Item($r, 3).Value = "s__F0040 sp900095835"
    $ws.Cells.Item($r, 4).Value = "s__F0040 sp900095835"
    $ws.Cells.Item($r, 5).Value = 0.5
}
